# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 96-97) before the existing data that
# currently starts at row 96, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 96-97; this shifts existing rows 96:111 down to 98:113
$ws.Range("A96:T97").EntireRow.Insert()

# --- New row 96 ---
$ws.Range("A96").Value = 9
$ws.Range("B96").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C96").Value = "Metropolitana"
$ws.Range("D96").Value = 44900
$ws.Range("E96").Value = 13
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100103
$ws.Range("H96").Value = "Frutos de hueso (carozo)"
$ws.Range("I96").Value = 100103003
$ws.Range("J96").Value = "Damasco"
$ws.Range("K96").Value = "Castle Brite"
$ws.Range("L96").Value = "Especial"
$ws.Range("M96").Value = 260
$ws.Range("N96").Value = 16000
$ws.Range("O96").Value = 16000
$ws.Range("P96").Value = 16000
$ws.Range("Q96").Value = "$/caja 16 kilos"
$ws.Range("R96").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S96").Value = 1000
$ws.Range("T96").Value = 16

# --- New row 97 ---
$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44900
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100103
$ws.Range("H97").Value = "Frutos de hueso (carozo)"
$ws.Range("I97").Value = 100103003
$ws.Range("J97").Value = "Damasco"
$ws.Range("K97").Value = "Castle Brite"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 300
$ws.Range("N97").Value = 14000
$ws.Range("O97").Value = 14000
$ws.Range("P97").Value = 14000
$ws.Range("Q97").Value = "$/caja 16 kilos"
$ws.Range("R97").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S97").Value = 800
$ws.Range("T97").Value = 16
